$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.049.99"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").Value = "2.289.37"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "505.68"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("D6").Value = "129.12"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("D9").Value = "2.307.62"
$ws.Range("E9").Value = "  +1.44%  "
$ws.Range("D10").Value = "0.0980"
$ws.Range("E10").Value = "  +3.09%  "
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").Value = "5.11"
$ws.Range("E12").Value = "  +8.60%  "
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("D14").Value = "23.71"
$ws.Range("E14").Value = "  +4.64%  "
$ws.Range("D15").Value = "2.697.50"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "55.115.44"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("D17").Value = "0.0000131"
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D18").Value = "2.307.44"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").Value = "10.44"
$ws.Range("E19").Value = "  +2.25%  "
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").Value = "313.05"
$ws.Range("E21").Value = "  +3.22%  "
$ws.Range("D22").Value = "6.63"
$ws.Range("E22").Value = "  +5.04%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").Value = "60.16"
$ws.Range("E24").Value = "  -1.53%  "
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("D27").Value = "7.50"
$ws.Range("E27").Value = "  +3.05%  "
$ws.Range("D28").Value = "172.28"
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0711"
$ws.Range("E29").Value = "  +4.07%  "
$ws.Range("E30").Value = "  +6.94%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "6.14"
$ws.Range("E31").Value = "  +3.96%  "
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "18.04"
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("D35").Value = "0.994"
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("E36").Value = "  +3.87%  "
$ws.Range("D37").Value = "0.916"
$ws.Range("E37").Value = "  -4.34%  "
$ws.Range("D38").Value = "3.89"
$ws.Range("E38").Value = "  +5.48%  "
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("E40").Value = "  +3.71%  "
$ws.Range("D41").Value = "0.376"
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("D42").Value = "136.12"
$ws.Range("E42").Value = "  +9.10%  "
$ws.Range("D43").Value = "5.13"
$ws.Range("E43").Value = "  +6.37%  "
$ws.Range("D44").Value = "3.44"
$ws.Range("D45").Value = "261.66"
$ws.Range("E45").Value = "  +9.78%  "
$ws.Range("E46").Value = "  +3.24%  "
$ws.Range("D47").Value = "0.0914"
$ws.Range("E47").Value = "  +2.49%  "
$ws.Range("D48").Value = "0.553"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("E49").Value = "  +1.13%  "
$ws.Range("D50").Value = "0.0211"
$ws.Range("E50").Value = "  +3.38%  "
$ws.Range("D51").Value = "16.50"
$ws.Range("E51").Value = "  +2.01%  "
